# Penalty / Reward System attempt (unfinished) - remove some weekly/monthly
# PO rows and tweak a couple of values, per the commit diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: "Weekly Quantity" ----
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

# Delete rows (original row numbers) 16 and 17 first (bottom-most first so
# the remaining row numbers we still need to delete stay valid), then 7,
# then 5 and 4.
$ws1.Rows.Item(17).Delete()
$ws1.Rows.Item(16).Delete()
$ws1.Rows.Item(7).Delete()
$ws1.Rows.Item(5).Delete()
$ws1.Rows.Item(4).Delete()

# ---- Sheet 2: "Monthly Trend" ----
$ws2 = $wb.Worksheets.Item("Monthly Trend")

# Update requested quantity values for rows 3 and 4.
$ws2.Cells.Item(3, 2).Value = 30
$ws2.Cells.Item(4, 2).Value = 150

# Delete row 8 (45382.99999999999 / 80), shifting subsequent rows up.
$ws2.Rows.Item(8).Delete()
